$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 46.595173
$ws.Cells.Item(2, 8).Value = 139.785519
$ws.Cells.Item(2, 9).Value = 0.7981698877785356
$ws.Cells.Item(2, 10).Value = 0.7981698877785355
$ws.Cells.Item(2, 13).Value = 0.1825283333333333
$ws.Cells.Item(2, 14).Value = 0.547585
$ws.Cells.Item(2, 15).Value = 0.001028331058213739
$ws.Cells.Item(2, 16).Value = 0.001028331058213739
$ws.Cells.Item(2, 17).Value = 8.504939269068332
$ws.Cells.Item(2, 18).Value = 76.54445342161499
$ws.Cells.Item(2, 19).Value = 0.0008207828853336426
$ws.Cells.Item(2, 20).Value = 0.0008207828853336426

# Row 3
$ws.Cells.Item(3, 7).Value = 46.595173
$ws.Cells.Item(3, 8).Value = 139.785519
$ws.Cells.Item(3, 9).Value = 0.7981698877785356
$ws.Cells.Item(3, 10).Value = 0.7981698877785355
$ws.Cells.Item(3, 15).Value = 0.0001759459539160193
$ws.Cells.Item(3, 16).Value = 0.0001759459539160193
$ws.Cells.Item(3, 17).Value = 1.455182784514333
$ws.Cells.Item(3, 18).Value = 13.096645060629
$ws.Cells.Item(3, 19).Value = 0.0001404347622922365
$ws.Cells.Item(3, 20).Value = 0.0001404347622922365

# Row 4
$ws.Cells.Item(4, 7).Value = 46.595173
$ws.Cells.Item(4, 8).Value = 139.785519
$ws.Cells.Item(4, 9).Value = 0.7981698877785356
$ws.Cells.Item(4, 10).Value = 0.7981698877785355
$ws.Cells.Item(4, 13).Value = 103.239782
$ws.Cells.Item(4, 14).Value = 309.719346
$ws.Cells.Item(4, 15).Value = 0.5816339432625932
$ws.Cells.Item(4, 16).Value = 0.5816339432625932
$ws.Cells.Item(4, 17).Value = 4810.475502772285
$ws.Cells.Item(4, 18).Value = 43294.27952495057
$ws.Cells.Item(4, 19).Value = 0.4642426992220912
$ws.Cells.Item(4, 20).Value = 0.4642426992220911

# Row 5
$ws.Cells.Item(5, 7).Value = 46.595173
$ws.Cells.Item(5, 8).Value = 139.785519
$ws.Cells.Item(5, 9).Value = 0.7981698877785356
$ws.Cells.Item(5, 10).Value = 0.7981698877785355
$ws.Cells.Item(5, 13).Value = 0.04852733333333333
$ws.Cells.Item(5, 14).Value = 0.145582
$ws.Cells.Item(5, 15).Value = 0.0002733940705404138
$ws.Cells.Item(5, 16).Value = 0.0002733940705404139
$ws.Cells.Item(5, 17).Value = 2.261139491895333
$ws.Cells.Item(5, 18).Value = 20.350255427058
$ws.Cells.Item(5, 19).Value = 0.0002182149146025592
$ws.Cells.Item(5, 20).Value = 0.0002182149146025592

# Row 6
$ws.Cells.Item(6, 7).Value = 46.595173
$ws.Cells.Item(6, 8).Value = 139.785519
$ws.Cells.Item(6, 9).Value = 0.7981698877785356
$ws.Cells.Item(6, 10).Value = 0.7981698877785355
$ws.Cells.Item(6, 13).Value = 73.99751433333334
$ws.Cells.Item(6, 14).Value = 221.992543
$ws.Cells.Item(6, 15).Value = 0.4168883856547366
$ws.Cells.Item(6, 16).Value = 0.4168883856547366
$ws.Cells.Item(6, 17).Value = 3447.926981931646
$ws.Cells.Item(6, 18).Value = 31031.34283738482
$ws.Cells.Item(6, 19).Value = 0.332747755994216
$ws.Cells.Item(6, 20).Value = 0.332747755994216

# Row 7
$ws.Cells.Item(7, 9).Value = 0.02837882113957134
$ws.Cells.Item(7, 10).Value = 0.02837882113957133
$ws.Cells.Item(7, 13).Value = 0.1825283333333333
$ws.Cells.Item(7, 14).Value = 0.547585
$ws.Cells.Item(7, 15).Value = 0.001028331058213739
$ws.Cells.Item(7, 16).Value = 0.001028331058213739
$ws.Cells.Item(7, 17).Value = 0.3023919519083333
$ws.Cells.Item(7, 18).Value = 2.721527567175
$ws.Cells.Item(7, 19).Value = 0.00002918282317331381
$ws.Cells.Item(7, 20).Value = 0.00002918282317331381

# Row 8
$ws.Cells.Item(8, 9).Value = 0.02837882113957134
$ws.Cells.Item(8, 10).Value = 0.02837882113957133
$ws.Cells.Item(8, 15).Value = 0.0001759459539160193
$ws.Cells.Item(8, 16).Value = 0.0001759459539160193
$ws.Cells.Item(8, 19).Value = 0.000004993138756413971
$ws.Cells.Item(8, 20).Value = 0.00000499313875641397

# Row 9
$ws.Cells.Item(9, 9).Value = 0.02837882113957134
$ws.Cells.Item(9, 10).Value = 0.02837882113957133
$ws.Cells.Item(9, 13).Value = 103.239782
$ws.Cells.Item(9, 14).Value = 309.719346
$ws.Cells.Item(9, 15).Value = 0.5816339432625932
$ws.Cells.Item(9, 16).Value = 0.5816339432625932
$ws.Cells.Item(9, 17).Value = 171.03579824267
$ws.Cells.Item(9, 18).Value = 1539.32218418403
$ws.Cells.Item(9, 19).Value = 0.01650608564455271
$ws.Cells.Item(9, 20).Value = 0.01650608564455271

# Row 10
$ws.Cells.Item(10, 9).Value = 0.02837882113957134
$ws.Cells.Item(10, 10).Value = 0.02837882113957133
$ws.Cells.Item(10, 13).Value = 0.04852733333333333
$ws.Cells.Item(10, 14).Value = 0.145582
$ws.Cells.Item(10, 15).Value = 0.0002733940705404138
$ws.Cells.Item(10, 16).Value = 0.0002733940705404139
$ws.Cells.Item(10, 17).Value = 0.08039450522333334
$ws.Cells.Item(10, 18).Value = 0.72355054701
$ws.Cells.Item(10, 19).Value = 0.000007758601428485753
$ws.Cells.Item(10, 20).Value = 0.000007758601428485753

# Row 11
$ws.Cells.Item(11, 9).Value = 0.02837882113957134
$ws.Cells.Item(11, 10).Value = 0.02837882113957133
$ws.Cells.Item(11, 13).Value = 73.99751433333334
$ws.Cells.Item(11, 14).Value = 221.992543
$ws.Cells.Item(11, 15).Value = 0.4168883856547366
$ws.Cells.Item(11, 16).Value = 0.4168883856547366
$ws.Cells.Item(11, 17).Value = 122.5905720333184
$ws.Cells.Item(11, 18).Value = 1103.315148299865
$ws.Cells.Item(11, 19).Value = 0.01183080093166041
$ws.Cells.Item(11, 20).Value = 0.0118308009316604

# Row 12
$ws.Cells.Item(12, 7).Value = 5.966798333333333
$ws.Cells.Item(12, 8).Value = 17.900395
$ws.Cells.Item(12, 9).Value = 0.102210560654294
$ws.Cells.Item(12, 10).Value = 0.1022105606542939
$ws.Cells.Item(12, 13).Value = 0.1825283333333333
$ws.Cells.Item(12, 14).Value = 0.547585
$ws.Cells.Item(12, 15).Value = 0.001028331058213739
$ws.Cells.Item(12, 16).Value = 0.001028331058213739
$ws.Cells.Item(12, 17).Value = 1.089109755119444
$ws.Cells.Item(12, 18).Value = 9.801987796075
$ws.Cells.Item(12, 19).Value = 0.0001051062939982496
$ws.Cells.Item(12, 20).Value = 0.0001051062939982496

# Row 13
$ws.Cells.Item(13, 7).Value = 5.966798333333333
$ws.Cells.Item(13, 8).Value = 17.900395
$ws.Cells.Item(13, 9).Value = 0.102210560654294
$ws.Cells.Item(13, 10).Value = 0.1022105606542939
$ws.Cells.Item(13, 15).Value = 0.0001759459539160193
$ws.Cells.Item(13, 16).Value = 0.0001759459539160193
$ws.Cells.Item(13, 17).Value = 0.1863451008827778
$ws.Cells.Item(13, 18).Value = 1.677105907945
$ws.Cells.Item(13, 19).Value = 0.0000179835345946109
$ws.Cells.Item(13, 20).Value = 0.00001798353459461089

# Row 14
$ws.Cells.Item(14, 7).Value = 5.966798333333333
$ws.Cells.Item(14, 8).Value = 17.900395
$ws.Cells.Item(14, 9).Value = 0.102210560654294
$ws.Cells.Item(14, 10).Value = 0.1022105606542939
$ws.Cells.Item(14, 13).Value = 103.239782
$ws.Cells.Item(14, 14).Value = 309.719346
$ws.Cells.Item(14, 15).Value = 0.5816339432625932
$ws.Cells.Item(14, 16).Value = 0.5816339432625932
$ws.Cells.Item(14, 17).Value = 616.0109591712966
$ws.Cells.Item(14, 18).Value = 5544.09863254167
$ws.Cells.Item(14, 19).Value = 0.05944913143643745
$ws.Cells.Item(14, 20).Value = 0.05944913143643744

# Row 15
$ws.Cells.Item(15, 7).Value = 5.966798333333333
$ws.Cells.Item(15, 8).Value = 17.900395
$ws.Cells.Item(15, 9).Value = 0.102210560654294
$ws.Cells.Item(15, 10).Value = 0.1022105606542939
$ws.Cells.Item(15, 13).Value = 0.04852733333333333
$ws.Cells.Item(15, 14).Value = 0.145582
$ws.Cells.Item(15, 15).Value = 0.0002733940705404138
$ws.Cells.Item(15, 16).Value = 0.0002733940705404139
$ws.Cells.Item(15, 17).Value = 0.2895528116544444
$ws.Cells.Item(15, 18).Value = 2.60597530489
$ws.Cells.Item(15, 19).Value = 0.00002794376122949529
$ws.Cells.Item(15, 20).Value = 0.00002794376122949529

# Row 16
$ws.Cells.Item(16, 7).Value = 5.966798333333333
$ws.Cells.Item(16, 8).Value = 17.900395
$ws.Cells.Item(16, 9).Value = 0.102210560654294
$ws.Cells.Item(16, 10).Value = 0.1022105606542939
$ws.Cells.Item(16, 13).Value = 73.99751433333334
$ws.Cells.Item(16, 14).Value = 221.992543
$ws.Cells.Item(16, 15).Value = 0.4168883856547366
$ws.Cells.Item(16, 16).Value = 0.4168883856547366
$ws.Cells.Item(16, 17).Value = 441.5282451949428
$ws.Cells.Item(16, 18).Value = 3973.754206754485
$ws.Cells.Item(16, 19).Value = 0.04261039562803415
$ws.Cells.Item(16, 20).Value = 0.04261039562803414

# Row 17
$ws.Cells.Item(17, 7).Value = 0.4842143333333333
$ws.Cells.Item(17, 8).Value = 1.452643
$ws.Cells.Item(17, 9).Value = 0.008294535146321381
$ws.Cells.Item(17, 10).Value = 0.008294535146321381
$ws.Cells.Item(17, 13).Value = 0.1825283333333333
$ws.Cells.Item(17, 14).Value = 0.547585
$ws.Cells.Item(17, 15).Value = 0.001028331058213739
$ws.Cells.Item(17, 16).Value = 0.001028331058213739
$ws.Cells.Item(17, 17).Value = 0.08838283523944443
$ws.Cells.Item(17, 18).Value = 0.7954455171549999
$ws.Cells.Item(17, 19).Value = 0.000008529528104407714
$ws.Cells.Item(17, 20).Value = 0.000008529528104407716

# Row 18
$ws.Cells.Item(18, 7).Value = 0.4842143333333333
$ws.Cells.Item(18, 8).Value = 1.452643
$ws.Cells.Item(18, 9).Value = 0.008294535146321381
$ws.Cells.Item(18, 10).Value = 0.008294535146321381
$ws.Cells.Item(18, 15).Value = 0.0001759459539160193
$ws.Cells.Item(18, 16).Value = 0.0001759459539160193
$ws.Cells.Item(18, 17).Value = 0.01512217503477778
$ws.Cells.Item(18, 18).Value = 0.136099575313
$ws.Cells.Item(18, 19).Value = 0.000001459389898609464
$ws.Cells.Item(18, 20).Value = 0.000001459389898609464

# Row 19
$ws.Cells.Item(19, 7).Value = 0.4842143333333333
$ws.Cells.Item(19, 8).Value = 1.452643
$ws.Cells.Item(19, 9).Value = 0.008294535146321381
$ws.Cells.Item(19, 10).Value = 0.008294535146321381
$ws.Cells.Item(19, 13).Value = 103.239782
$ws.Cells.Item(19, 14).Value = 309.719346
$ws.Cells.Item(19, 15).Value = 0.5816339432625932
$ws.Cells.Item(19, 16).Value = 0.5816339432625932
$ws.Cells.Item(19, 17).Value = 49.99018221460866
$ws.Cells.Item(19, 18).Value = 449.9116399314779
$ws.Cells.Item(19, 19).Value = 0.004824383184685075
$ws.Cells.Item(19, 20).Value = 0.004824383184685075

# Row 20
$ws.Cells.Item(20, 7).Value = 0.4842143333333333
$ws.Cells.Item(20, 8).Value = 1.452643
$ws.Cells.Item(20, 9).Value = 0.008294535146321381
$ws.Cells.Item(20, 10).Value = 0.008294535146321381
$ws.Cells.Item(20, 13).Value = 0.04852733333333333
$ws.Cells.Item(20, 14).Value = 0.145582
$ws.Cells.Item(20, 15).Value = 0.0002733940705404138
$ws.Cells.Item(20, 16).Value = 0.0002733940705404139
$ws.Cells.Item(20, 17).Value = 0.02349763035844444
$ws.Cells.Item(20, 18).Value = 0.211478673226
$ws.Cells.Item(20, 19).Value = 0.000002267676726893329
$ws.Cells.Item(20, 20).Value = 0.00000226767672689333

# Row 21
$ws.Cells.Item(21, 7).Value = 0.4842143333333333
$ws.Cells.Item(21, 8).Value = 1.452643
$ws.Cells.Item(21, 9).Value = 0.008294535146321381
$ws.Cells.Item(21, 10).Value = 0.008294535146321381
$ws.Cells.Item(21, 13).Value = 73.99751433333334
$ws.Cells.Item(21, 14).Value = 221.992543
$ws.Cells.Item(21, 15).Value = 0.4168883856547366
$ws.Cells.Item(21, 16).Value = 0.4168883856547366
$ws.Cells.Item(21, 17).Value = 35.83065707123878
$ws.Cells.Item(21, 18).Value = 322.475913641149
$ws.Cells.Item(21, 19).Value = 0.003457895366906395
$ws.Cells.Item(21, 20).Value = 0.003457895366906395

# Row 22
$ws.Cells.Item(22, 7).Value = 3.674642333333333
$ws.Cells.Item(22, 8).Value = 11.023927
$ws.Cells.Item(22, 9).Value = 0.0629461952812778
$ws.Cells.Item(22, 10).Value = 0.0629461952812778
$ws.Cells.Item(22, 13).Value = 0.1825283333333333
$ws.Cells.Item(22, 14).Value = 0.547585
$ws.Cells.Item(22, 15).Value = 0.001028331058213739
$ws.Cells.Item(22, 16).Value = 0.001028331058213739
$ws.Cells.Item(22, 17).Value = 0.6707263406994444
$ws.Cells.Item(22, 18).Value = 6.036537066295
$ws.Cells.Item(22, 19).Value = 0.00006472952760412504
$ws.Cells.Item(22, 20).Value = 0.00006472952760412506

# Row 23
$ws.Cells.Item(23, 7).Value = 3.674642333333333
$ws.Cells.Item(23, 8).Value = 11.023927
$ws.Cells.Item(23, 9).Value = 0.0629461952812778
$ws.Cells.Item(23, 10).Value = 0.0629461952812778
$ws.Cells.Item(23, 15).Value = 0.0001759459539160193
$ws.Cells.Item(23, 16).Value = 0.0001759459539160193
$ws.Cells.Item(23, 17).Value = 0.1147603049507778
$ws.Cells.Item(23, 18).Value = 1.032842744557
$ws.Cells.Item(23, 19).Value = 0.00001107512837414845
$ws.Cells.Item(23, 20).Value = 0.00001107512837414845

# Row 24
$ws.Cells.Item(24, 7).Value = 3.674642333333333
$ws.Cells.Item(24, 8).Value = 11.023927
$ws.Cells.Item(24, 9).Value = 0.0629461952812778
$ws.Cells.Item(24, 10).Value = 0.0629461952812778
$ws.Cells.Item(24, 13).Value = 103.239782
$ws.Cells.Item(24, 14).Value = 309.719346
$ws.Cells.Item(24, 15).Value = 0.5816339432625932
$ws.Cells.Item(24, 16).Value = 0.5816339432625932
$ws.Cells.Item(24, 17).Value = 379.3692734213046
$ws.Cells.Item(24, 18).Value = 3414.323460791742
$ws.Cells.Item(24, 19).Value = 0.03661164377482685
$ws.Cells.Item(24, 20).Value = 0.03661164377482685

# Row 25
$ws.Cells.Item(25, 7).Value = 3.674642333333333
$ws.Cells.Item(25, 8).Value = 11.023927
$ws.Cells.Item(25, 9).Value = 0.0629461952812778
$ws.Cells.Item(25, 10).Value = 0.0629461952812778
$ws.Cells.Item(25, 13).Value = 0.04852733333333333
$ws.Cells.Item(25, 14).Value = 0.145582
$ws.Cells.Item(25, 15).Value = 0.0002733940705404138
$ws.Cells.Item(25, 16).Value = 0.0002733940705404139
$ws.Cells.Item(25, 17).Value = 0.1783205933904445
$ws.Cells.Item(25, 18).Value = 1.604885340514
$ws.Cells.Item(25, 19).Value = 0.00001720911655298033
$ws.Cells.Item(25, 20).Value = 0.00001720911655298033

# Row 26
$ws.Cells.Item(26, 7).Value = 3.674642333333333
$ws.Cells.Item(26, 8).Value = 11.023927
$ws.Cells.Item(26, 9).Value = 0.0629461952812778
$ws.Cells.Item(26, 10).Value = 0.0629461952812778
$ws.Cells.Item(26, 13).Value = 73.99751433333334
$ws.Cells.Item(26, 14).Value = 221.992543
$ws.Cells.Item(26, 15).Value = 0.4168883856547366
$ws.Cells.Item(26, 16).Value = 0.4168883856547366
$ws.Cells.Item(26, 17).Value = 271.9143987307068
$ws.Cells.Item(26, 18).Value = 2447.229588576361
$ws.Cells.Item(26, 19).Value = 0.02624153773391971
$ws.Cells.Item(26, 20).Value = 0.02624153773391971
